# Apply cryptos list update (price/volume refresh) as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to be written back as plain text (matching the source
    # inline-string cells), even when the text looks like a number, and then
    # restore the default "Normal" style so no stray formatting is introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '41.921.11'
Set-TextValue $ws.Range('E2') '  +0.99%  '
Set-TextValue $ws.Range('D3') '2.143.73'
Set-TextValue $ws.Range('E3') '  -1.01%  '
Set-TextValue $ws.Range('E4') '  -0.14%  '
Set-TextValue $ws.Range('D5') '251.85'
Set-TextValue $ws.Range('E5') '  +5.88%  '
Set-TextValue $ws.Range('D6') '0.604'
Set-TextValue $ws.Range('E6') '  -0.54%  '
Set-TextValue $ws.Range('D7') '72.39'
Set-TextValue $ws.Range('E7') '  +0.90%  '
Set-TextValue $ws.Range('E8') '  -0.07%  '
Set-TextValue $ws.Range('D9') '0.573'
Set-TextValue $ws.Range('E9') '  -0.60%  '
Set-TextValue $ws.Range('D10') '39.29'
Set-TextValue $ws.Range('E10') '  -1.23%  '
Set-TextValue $ws.Range('D11') '0.0902'
Set-TextValue $ws.Range('E11') '  -0.48%  '
Set-TextValue $ws.Range('E12') '  +0.53%  '
Set-TextValue $ws.Range('D13') '6.67'
Set-TextValue $ws.Range('E13') '  -0.35%  '
Set-TextValue $ws.Range('D14') '2.466.76'
Set-TextValue $ws.Range('E14') '  -1.02%  '
Set-TextValue $ws.Range('D15') '14.04'
Set-TextValue $ws.Range('E15') '  -1.60%  '
Set-TextValue $ws.Range('D16') '2.148.78'
Set-TextValue $ws.Range('E16') '  -0.88%  '
Set-TextValue $ws.Range('D17') '0.760'
Set-TextValue $ws.Range('E17') '  -2.27%  '
Set-TextValue $ws.Range('D18') '41.788.46'
Set-TextValue $ws.Range('E18') '  +0.93%  '
Set-TextValue $ws.Range('D19') '0.0000101'
Set-TextValue $ws.Range('E19') '  -1.69%  '
Set-TextValue $ws.Range('D20') '70.00'
Set-TextValue $ws.Range('E20') '  +0.02%  '
Set-TextValue $ws.Range('D21') '5.78'
Set-TextValue $ws.Range('E21') '  +0.22%  '
Set-TextValue $ws.Range('B22') 'BitcoinCash'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D22') '224.70'
Set-TextValue $ws.Range('E22') '  -0.66%  '
Set-TextValue $ws.Range('B23') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D23') '9.48'
Set-TextValue $ws.Range('E23') '  -1.96%  '
Set-TextValue $ws.Range('D24') '2.12'
Set-TextValue $ws.Range('E24') '  +5.16%  '
Set-TextValue $ws.Range('E25') '  -0.20%  '
Set-TextValue $ws.Range('D26') '10.36'
Set-TextValue $ws.Range('E26') '  -3.11%  '
Set-TextValue $ws.Range('E27') '  +0.72%  '
Set-TextValue $ws.Range('E28') '  +2.93%  '
Set-TextValue $ws.Range('D29') '2.18'
Set-TextValue $ws.Range('E29') '  -0.61%  '
Set-TextValue $ws.Range('D30') '36.33'
Set-TextValue $ws.Range('E30') '  +8.46%  '
Set-TextValue $ws.Range('D31') '167.83'
Set-TextValue $ws.Range('E31') '  -1.83%  '
Set-TextValue $ws.Range('D32') '19.79'
Set-TextValue $ws.Range('E32') '  +0.04%  '
Set-TextValue $ws.Range('D33') '0.0791'
Set-TextValue $ws.Range('E33') '  +3.00%  '
Set-TextValue $ws.Range('D34') '5.06'
Set-TextValue $ws.Range('E34') '  -2.81%  '
Set-TextValue $ws.Range('D35') '0.119'
Set-TextValue $ws.Range('D36') '0.105'
Set-TextValue $ws.Range('E36') '  +1.34%  '
Set-TextValue $ws.Range('D37') '4.20'
Set-TextValue $ws.Range('E37') '  -2.42%  '
Set-TextValue $ws.Range('D38') '0.0324'
Set-TextValue $ws.Range('E38') '  +6.52%  '
Set-TextValue $ws.Range('D39') '11.77'
Set-TextValue $ws.Range('E39') '  -2.88%  '
Set-TextValue $ws.Range('D40') '2.03'
Set-TextValue $ws.Range('E40') '  -3.27%  '
Set-TextValue $ws.Range('D41') '0.193'
Set-TextValue $ws.Range('E41') '  +3.11%  '
Set-TextValue $ws.Range('D42') '5.08'
Set-TextValue $ws.Range('E42') '  -5.33%  '
Set-TextValue $ws.Range('D43') '57.95'
Set-TextValue $ws.Range('E43') '  -1.26%  '
Set-TextValue $ws.Range('D44') '99.39'
Set-TextValue $ws.Range('E44') '  +2.28%  '
Set-TextValue $ws.Range('D45') '0.461'
Set-TextValue $ws.Range('E45') '  +13.67%  '
Set-TextValue $ws.Range('D46') '8.19'
Set-TextValue $ws.Range('E46') '  -3.08%  '
Set-TextValue $ws.Range('D47') '0.0956'
Set-TextValue $ws.Range('E47') '  -0.78%  '
Set-TextValue $ws.Range('D48') '2.35'
Set-TextValue $ws.Range('E48') '  +7.77%  '
Set-TextValue $ws.Range('D49') '1.07'
Set-TextValue $ws.Range('E49') '  -0.33%  '
Set-TextValue $ws.Range('E50') '  +0.33%  '
Set-TextValue $ws.Range('E51') '  +0.68%  '
